$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.574.19'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '3.586.43'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.78'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.32'
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = '3.586.29'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.136'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.99'
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.415'
$ws.Range("D13").Value = '4.199.49'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.00'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").Value = '3.592.68'
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("D17").Value = '66.685.16'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.47'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.02'
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.50'
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.623'
$ws.Range("E23").Value = '  +2.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.10'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = '3.732.17'
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000121'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.30'
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.08'
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = '3.585.93'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.47'
$ws.Range("E33").Value = '  +0.21%  '
$ws.Range("E34").Value = '  -3.14%  '
$ws.Range("E35").Value = '  -1.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.84'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.63'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.13'
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0853'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.22'
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.894'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '45.86'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.53'
$ws.Range("E47").Value = '  +5.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.19'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.93'
$ws.Range("E49").Value = '  -4.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.72'
$ws.Range("E50").Value = '  +4.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.20'
$ws.Range("E51").Value = '  +1.00%  '
